# Apply a row-wise shuffle of the price/date/quality columns (D, I, J, K, L, M, N, O, P, Q)
# across data rows 2..111. Each destination row takes its new values from the source row
# indicated by $mapping (1-based offset into data rows, i.e. $mapping[0] is the source row
# for row 2, $mapping[1] is the source row for row 3, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 111

# For each destination row (index 0 => row 2), the source row number to copy the
# "moving" columns from.
$mapping = @(42,92,48,19,88,53,28,10,25,56,70,35,95,60,54,44,58,41,61,96,9,80,87,85,46,47,99,7,8,67,57,55,108,31,15,69,63,110,23,2,13,4,36,83,22,45,16,29,90,91,103,5,105,97,106,14,109,17,94,20,82,37,39,64,26,65,101,18,73,40,3,62,59,98,77,71,11,107,34,6,51,12,104,38,86,33,102,66,24,111,79,50,72,100,75,76,78,32,81,30,27,89,74,93,49,84,21,52,43,68)

# Column letters for the values that move together as a row-unit.
$cols = @("D","I","J","K","L","M","N","O","P","Q")

# Snapshot all the "moving" column values for every data row BEFORE any writes happen,
# so that source rows read later in the loop are never polluted by earlier writes.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the shuffled values back into each destination row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r - $firstRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}
